$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 339, shifting existing rows 339:451 down to 340:452.
$ws.Rows.Item(339).Insert()

# Populate the newly inserted row 339 with the new record's data.
$ws.Range("A339").Value2 = 1
$ws.Range("B339").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C339").Value2 = "Arica y Parinacota"
$ws.Range("D339").Value2 = 44985
$ws.Range("E339").Value2 = 15
$ws.Range("F339").Value2 = 100114013
$ws.Range("G339").Value2 = "Zanahoria"
$ws.Range("H339").Value2 = "Sin especificar"
$ws.Range("I339").Value2 = "Primera"
$ws.Range("J339").Value2 = 15
$ws.Range("K339").Value2 = 9000
$ws.Range("L339").Value2 = 10000
$ws.Range("M339").Value2 = 9667
$ws.Range("N339").Value2 = "$/saco 25 kilos"
$ws.Range("O339").Value2 = "Valle de Camiña"
$ws.Range("P339").Value2 = 387
$ws.Range("Q339").Value2 = 25
$ws.Range("R339").Value2 = "Hortaliza"
